$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.270.94"
$ws.Range("E2").Value = "  +0.08%  "
$ws.Range("D3").Value = "3.417.15"
$ws.Range("E3").Value = "  +1.08%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "578.97"
$ws.Range("E5").Value = "  -1.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "175.83"
$ws.Range("E6").Value = "  -2.79%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").Value = "3.409.80"
$ws.Range("E8").Value = "  +0.93%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("E10").Value = "  +0.74%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.582"
$ws.Range("E11").Value = "  -0.90%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "48.44"
$ws.Range("E12").Value = "  -0.98%  "
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "695.80"
$ws.Range("E14").Value = "  +1.34%  "
$ws.Range("D15").Value = "3.952.37"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("E16").Value = "  +0.11%  "
$ws.Range("D17").Value = "69.290.77"
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "3.424.08"
$ws.Range("E18").Value = "  +1.71%  "
$ws.Range("E19").Value = "  +0.83%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "17.65"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.32"
$ws.Range("E21").Value = "  -0.68%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.894"
$ws.Range("E22").Value = "  -0.69%  "
$ws.Range("E23").Value = "  +0.15%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "16.87"
$ws.Range("E24").Value = "  -1.38%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "100.24"
$ws.Range("E25").Value = "  -4.31%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "3.87"
$ws.Range("E26").Value = "  -2.02%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "2.65"
$ws.Range("E27").Value = "  -2.68%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.55"
$ws.Range("E28").Value = "  -0.56%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "33.27"
$ws.Range("E29").Value = "  -3.40%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.69"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.89"
$ws.Range("E31").Value = "  -1.84%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "565.83"
$ws.Range("E32").Value = "  +1.49%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.71"
$ws.Range("E33").Value = "  +0.91%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "10.95"
$ws.Range("E34").Value = "  -2.06%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.104"
$ws.Range("E35").Value = "  -1.97%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "58.24"
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  -0.11%  "
$ws.Range("D38").Value = "3.579.68"
$ws.Range("E38").Value = "  -3.69%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.137"
$ws.Range("E39").Value = "  -3.09%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "34.67"
$ws.Range("E40").Value = "  -0.90%  "
$ws.Range("D41").Value = "0.0₃0723"
$ws.Range("E41").Value = "  +2.97%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.24"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.64"
$ws.Range("E43").Value = "  -0.90%  "
$ws.Range("B44").Value = "TheGraph"
$ws.Range("C44").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.331"
$ws.Range("E44").Value = "  -2.70%  "
$ws.Range("B45").Value = "VeChain"
$ws.Range("C45").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0416"
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("B46").Value = "Mantle"
$ws.Range("C46").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.45"
$ws.Range("E46").Value = "  +2.70%  "
$ws.Range("B47").Value = "ThetaToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.63"
$ws.Range("E47").Value = "  -1.18%  "
$ws.Range("B48").Value = "Stellar"
$ws.Range("C48").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.128"
$ws.Range("E48").Value = "  -1.53%  "
$ws.Range("B49").Value = "FirstDigitalUSD"
$ws.Range("C49").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.00"
$ws.Range("E49").Value = "  -0.11%  "
$ws.Range("B50").Value = "Monero"
$ws.Range("C50").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "131.65"
$ws.Range("E50").Value = "  -0.55%  "
$ws.Range("B51").Value = "CoreDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/HFvoXUQh4+coredao-core"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.63"
$ws.Range("E51").Value = "  +0.63%  "
